# Add new columns I0 (column I) and IF (column J) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Make the new header cells match the style/formatting of the existing header cells (e.g. H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2 through 32
$iValues = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 7, 3, 1, 4, 3, 2)
$jValues = @(6, 5, 5, 4, 3, 5, 5, 4, 5, 4, 6, 6, 6, 5, 5, 5, 5, 6, 5, 2, 6, 5, 6, 5, 4, 7, 4, 4, 6, 4, 2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
